$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.94943145514307
$ws.Range("C2").Value = 9.864324008596844
$ws.Range("D2").Value = 7.119280444993101
$ws.Range("E2").Value = 9.646071058415544
$ws.Range("F2").Value = 37.38304194697209
$ws.Range("I2").Value = 29.55346602580157
$ws.Range("L2").Value = 10.27488303849295
$ws.Range("M2").Value = 17.15870100078373
$ws.Range("N2").Value = 19.96945516764769
$ws.Range("B3").Value = 19.53062189691046
$ws.Range("C3").Value = 9.27920543866067
$ws.Range("D3").Value = 7.146797519895694
$ws.Range("E3").Value = 9.62575764367174
$ws.Range("F3").Value = 37.18173560654614
$ws.Range("I3").Value = 29.56658138265703
$ws.Range("L3").Value = 10.28598448907631
$ws.Range("M3").Value = 17.08375225915517
$ws.Range("N3").Value = 20.03970510188956
$ws.Range("B4").Value = 19.27467497848491
$ws.Range("C4").Value = 8.902594012468047
$ws.Range("D4").Value = 7.164480791490373
$ws.Range("E4").Value = 9.613020768527276
$ws.Range("F4").Value = 37.0690459123824
$ws.Range("I4").Value = 29.58179532446443
$ws.Range("L4").Value = 10.29436912991403
$ws.Range("M4").Value = 17.04140035155665
$ws.Range("N4").Value = 20.08481043669389
$ws.Range("B5").Value = 19.17083645641032
$ws.Range("C5").Value = 8.744893272948454
$ws.Range("D5").Value = 7.17188544495409
$ws.Range("E5").Value = 9.60776382062998
$ws.Range("F5").Value = 37.02589895041046
$ws.Range("I5").Value = 29.58979024439661
$ws.Range("L5").Value = 10.29818027449668
$ws.Range("M5").Value = 17.02507607616603
$ws.Range("N5").Value = 20.10368820845626
$ws.Range("B6").Value = 19.15362674653158
$ws.Range("C6").Value = 8.718456267096723
$ws.Range("D6").Value = 7.173126990610644
$ws.Range("E6").Value = 9.6068868850844
$ws.Range("F6").Value = 37.01890283923642
$ws.Range("I6").Value = 29.59122603435761
$ws.Range("L6").Value = 10.29883692684895
$ws.Range("M6").Value = 17.02242222061773
$ws.Range("N6").Value = 20.10685290020161
$ws.Range("B7").Value = 19.27327249349807
$ws.Range("C7").Value = 8.900484128665422
$ws.Range("D7").Value = 7.164579848557423
$ws.Range("E7").Value = 9.612950141121274
$ws.Range("F7").Value = 37.0684527445674
$ws.Range("I7").Value = 29.58189588660935
$ws.Range("L7").Value = 10.29441893178907
$ws.Range("M7").Value = 17.04117639761879
$ws.Range("N7").Value = 20.08506301512499
$ws.Range("B8").Value = 19.80487072862381
$ws.Range("C8").Value = 9.666246758342075
$ws.Range("D8").Value = 7.128605157179613
$ws.Range("E8").Value = 9.639121801429081
$ws.Range("F8").Value = 37.31138639159678
$ws.Range("I8").Value = 29.55649881818186
$ws.Range("L8").Value = 10.27838531863773
$ws.Range("M8").Value = 17.13210540160853
$ws.Range("N8").Value = 19.99326887228451
$ws.Range("B9").Value = 20.85005655740126
$ws.Range("C9").Value = 11.04298256727103
$ws.Range("D9").Value = 7.064285177217896
$ws.Range("E9").Value = 9.688352488207929
$ws.Range("F9").Value = 37.87279638311713
$ws.Range("I9").Value = 29.56372884662277
$ws.Range("L9").Value = 10.25938689211172
$ws.Range("M9").Value = 17.33891387664078
$ws.Range("N9").Value = 19.8288479805319
$ws.Range("B10").Value = 21.61080840404811
$ws.Range("C10").Value = 12.0250550814788
$ws.Range("D10").Value = 7.02079152767724
$ws.Range("E10").Value = 9.723257075612757
$ws.Range("F10").Value = 38.33468505712651
$ws.Range("I10").Value = 29.60404877431929
$ws.Range("L10").Value = 10.25301170929262
$ws.Range("M10").Value = 17.5073530962333
$ws.Range("N10").Value = 19.71746922808123
$ws.Range("B11").Value = 21.95360272834472
$ws.Range("C11").Value = 12.44626122269941
$ws.Range("D11").Value = 7.001815042190904
$ws.Range("E11").Value = 9.738863542812714
$ws.Range("F11").Value = 38.55494336143953
$ws.Range("I11").Value = 29.63002401485418
$ws.Range("L11").Value = 10.25175595347923
$ws.Range("M11").Value = 17.58735445500419
$ws.Range("N11").Value = 19.66882873274592
$ws.Range("B12").Value = 22.08280422908976
$ws.Range("C12").Value = 12.60210883532555
$ws.Range("D12").Value = 6.994744973427171
$ws.Range("E12").Value = 9.744734325829183
$ws.Range("F12").Value = 38.63975215803169
$ws.Range("I12").Value = 29.64095854831424
$ws.Range("L12").Value = 10.25151654102814
$ws.Range("M12").Value = 17.61811555998548
$ws.Range("N12").Value = 19.65069999959357
$ws.Range("B13").Value = 22.05500746367304
$ws.Range("C13").Value = 12.5687064726007
$ws.Range("D13").Value = 6.996262490435844
$ws.Range("E13").Value = 9.743471683755784
$ws.Range("F13").Value = 38.6214256593414
$ws.Range("I13").Value = 29.63855475674471
$ws.Range("L13").Value = 10.25155760782636
$ws.Range("M13").Value = 17.61147018239847
$ws.Range("N13").Value = 19.65459144828034
$ws.Range("B14").Value = 21.96424513507464
$ws.Range("C14").Value = 12.45915602205467
$ws.Range("D14").Value = 7.001231062635973
$ws.Range("E14").Value = 9.739347318028635
$ws.Range("F14").Value = 38.56189287504457
$ws.Range("I14").Value = 29.63090159487109
$ws.Range("L14").Value = 10.25173152702711
$ws.Range("M14").Value = 17.58987595749145
$ws.Range("N14").Value = 19.66733145726902
$ws.Range("B15").Value = 21.90856760135152
$ws.Range("C15").Value = 12.39157791446577
$ws.Range("D15").Value = 7.004289539839645
$ws.Range("E15").Value = 9.73681594086622
$ws.Range("F15").Value = 38.52560813789246
$ws.Range("I15").Value = 29.6263568392289
$ws.Range("L15").Value = 10.25186879515734
$ws.Range("M15").Value = 17.57670901235891
$ws.Range("N15").Value = 19.67517286544622
$ws.Range("B16").Value = 21.588328241366
$ws.Range("C16").Value = 11.99701590326105
$ws.Range("D16").Value = 7.022047928678736
$ws.Range("E16").Value = 9.722231690194908
$ws.Range("F16").Value = 38.32048997809468
$ws.Range("I16").Value = 29.60250491093561
$ws.Range("L16").Value = 10.25312683722667
$ws.Range("M16").Value = 17.50219114784106
$ws.Range("N16").Value = 19.72068868578451
$ws.Range("B17").Value = 21.3909332818912
$ws.Range("C17").Value = 11.74843618757459
$ws.Range("D17").Value = 7.033149026883205
$ws.Range("E17").Value = 9.713215137276755
$ws.Range("F17").Value = 38.19721508169562
$ws.Range("I17").Value = 29.58982839921137
$ws.Range("L17").Value = 10.25431954699092
$ws.Range("M17").Value = 17.4573284089087
$ws.Range("N17").Value = 19.74912941739427
$ws.Range("B18").Value = 21.27709625699106
$ws.Range("C18").Value = 11.60305428601889
$ws.Range("D18").Value = 7.039610255608824
$ws.Range("E18").Value = 9.708003468975587
$ws.Range("F18").Value = 38.12726892600386
$ws.Range("I18").Value = 29.58325573837204
$ws.Range("L18").Value = 10.25516036156638
$ws.Range("M18").Value = 17.43184366983731
$ws.Range("N18").Value = 19.76567855386962
$ws.Range("B19").Value = 21.23850564323904
$ws.Range("C19").Value = 11.55341675516556
$ws.Range("D19").Value = 7.04181101191471
$ws.Range("E19").Value = 9.706234492092891
$ws.Range("F19").Value = 38.1037526264918
$ws.Range("I19").Value = 29.58115372195608
$ws.Range("L19").Value = 10.2554716428437
$ws.Range("M19").Value = 17.42327036058393
$ws.Range("N19").Value = 19.77131460473426
$ws.Range("B20").Value = 21.41197844328641
$ws.Range("C20").Value = 11.77514673616293
$ws.Range("D20").Value = 7.031959415721316
$ws.Range("E20").Value = 9.714177612445738
$ws.Range("F20").Value = 38.21023912526195
$ws.Range("I20").Value = 29.59110346248006
$ws.Range("L20").Value = 10.25417656204047
$ws.Range("M20").Value = 17.46207123375736
$ws.Range("N20").Value = 19.74608211533773
$ws.Range("B21").Value = 21.99092175689953
$ws.Range("C21").Value = 12.49143266157362
$ws.Range("D21").Value = 6.999768529929384
$ws.Range("E21").Value = 9.740559803673914
$ws.Range("F21").Value = 38.57934152134642
$ws.Range("I21").Value = 29.63311970775713
$ws.Range("L21").Value = 10.25167403778122
$ws.Range("M21").Value = 17.59620620709453
$ws.Range("N21").Value = 19.66358153524612
$ws.Range("B22").Value = 22.36569777553016
$ws.Range("C22").Value = 12.93827738187661
$ws.Range("D22").Value = 6.979405377324579
$ws.Range("E22").Value = 9.757574507535047
$ws.Range("F22").Value = 38.82871335967636
$ws.Range("I22").Value = 29.66698085763254
$ws.Range("L22").Value = 10.25141453438906
$ws.Range("M22").Value = 17.68657974582431
$ws.Range("N22").Value = 19.61135466515733
$ws.Range("B23").Value = 22.16604472968813
$ws.Range("C23").Value = 12.70172963422862
$ws.Range("D23").Value = 6.990211906318963
$ws.Range("E23").Value = 9.748514252249789
$ws.Range("F23").Value = 38.69489344855786
$ws.Range("I23").Value = 29.64832290248318
$ws.Range("L23").Value = 10.25142726280296
$ws.Range("M23").Value = 17.63810445112077
$ws.Range("N23").Value = 19.63907463481921
$ws.Range("B24").Value = 21.40246501189593
$ws.Range("C24").Value = 11.76307859244116
$ws.Range("D24").Value = 7.032496992643375
$ws.Range("E24").Value = 9.713742564426102
$ws.Range("F24").Value = 38.20434806591594
$ws.Range("I24").Value = 29.59052477863685
$ws.Range("L24").Value = 10.25424072237921
$ws.Range("M24").Value = 17.45992604355631
$ws.Range("N24").Value = 19.74745918323156
$ws.Range("B25").Value = 20.56794662244372
$ws.Range("C25").Value = 10.67371753636316
$ws.Range("D25").Value = 7.081022366869353
$ws.Range("E25").Value = 9.675257170029523
$ws.Range("F25").Value = 37.71205026245573
$ws.Range("I25").Value = 29.55563854795075
$ws.Range("L25").Value = 10.2631943007378
$ws.Range("M25").Value = 17.2800019888988
$ws.Range("N25").Value = 19.87166793027195
